# Reproduces the "Add new Excel reader class" commit:
#   - Sheet2!A3 gets a new value "Hello" (new shared string, sheet grows
#     from an empty A1 dimension to A3).
#   - Sheet2's row-outline high-water mark (sheetFormatPr/@outlineLevelRow)
#     ends up at 2, with row 3 itself carrying no explicit outline level.
#   - Sheet2's remembered selection becomes A3 (A1 before), while Sheet1
#     stays the active/selected tab.
#   - Sheet1's remembered selection moves from I8 to I9.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: add the new cell value ---------------------------------------
$ws2.Range("A3").Value = "Hello"

# Reproduce the sheet-wide outlineLevelRow="2" metadata without leaving any
# grouped/outlined row behind in the used range: group a pair of far-away
# rows twice (outline level 1, then 2) which raises the sheet's outline
# high-water mark, then delete those helper rows again so Sheet2's real
# data (just row 3) is untouched.
$farRows = $ws2.Range("A1000:A1001")
$farRows.Rows.Group() | Out-Null
$farRows.Rows.Group() | Out-Null
$farRows.EntireRow.Delete() | Out-Null

# Leave Sheet2's remembered cursor on A3.
$ws2.Range("A3").Select() | Out-Null

# --- Sheet1: stays the active tab, cursor moves from I8 to I9 -------------
$ws1.Activate() | Out-Null
$ws1.Range("I9").Select() | Out-Null
